{"js": "// Update each three-digit\u00f7one-digit division answer cell with its new\n// problem/answer text. Each old value is unique in the document, so we\n// can safely search for the exact old text and replace it in place \u2014\n// this preserves the run's formatting (font, size, etc.) untouched.\nconst replacements = [\n  [\"519\u00f75=103, 4\", \"816\u00f78=102, 0\"],\n  [\"248\u00f76=41, 2\", \"778\u00f79=86, 4\"],\n  [\"441\u00f77=63, 0\", \"935\u00f78=116, 7\"],\n  [\"376\u00f75=75, 1\", \"219\u00f76=36, 3\"],\n  [\"381\u00f76=63, 3\", \"983\u00f76=163, 5\"],\n  [\"467\u00f79=51, 8\", \"563\u00f72=281, 1\"],\n  [\"609\u00f73=203, 0\", \"943\u00f72=471, 1\"],\n  [\"373\u00f79=41, 4\", \"554\u00f75=110, 4\"],\n  [\"170\u00f79=18, 8\", \"356\u00f76=59, 2\"],\n  [\"965\u00f79=107, 2\", \"981\u00f74=245, 1\"],\n  [\"691\u00f74=172, 3\", \"515\u00f78=64, 3\"],\n  [\"703\u00f75=140, 3\", \"159\u00f74=39, 3\"],\n  [\"413\u00f76=68, 5\", \"420\u00f79=46, 6\"],\n  [\"630\u00f76=105, 0\", \"795\u00f77=113, 4\"],\n  [\"875\u00f75=175, 0\", \"640\u00f75=128, 0\"],\n  [\"628\u00f75=125, 3\", \"367\u00f72=183, 1\"],\n  [\"601\u00f74=150, 1\", \"561\u00f78=70, 1\"],\n  [\"382\u00f78=47, 6\", \"554\u00f78=69, 2\"],\n  [\"745\u00f78=93, 1\", \"602\u00f76=100, 2\"],\n  [\"872\u00f74=218, 0\", \"544\u00f79=60, 4\"],\n  [\"348\u00f77=49, 5\", \"250\u00f75=50, 0\"],\n  [\"477\u00f78=59, 5\", \"546\u00f75=109, 1\"],\n  [\"436\u00f73=145, 1\", \"729\u00f74=182, 1\"],\n  [\"688\u00f72=344, 0\", \"744\u00f78=93, 0\"],\n  [\"896\u00f74=224, 0\", \"841\u00f76=140, 1\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update each three-digit\u00f7one-digit division answer cell with its new\n# problem/answer text. Each old value is unique in the document, so a\n# plain Find/Replace (wdReplaceOne) on the exact old string safely\n# targets the single matching run without touching its formatting.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"519\u00f75=103, 4\", \"816\u00f78=102, 0\"),\n    @(\"248\u00f76=41, 2\", \"778\u00f79=86, 4\"),\n    @(\"441\u00f77=63, 0\", \"935\u00f78=116, 7\"),\n    @(\"376\u00f75=75, 1\", \"219\u00f76=36, 3\"),\n    @(\"381\u00f76=63, 3\", \"983\u00f76=163, 5\"),\n    @(\"467\u00f79=51, 8\", \"563\u00f72=281, 1\"),\n    @(\"609\u00f73=203, 0\", \"943\u00f72=471, 1\"),\n    @(\"373\u00f79=41, 4\", \"554\u00f75=110, 4\"),\n    @(\"170\u00f79=18, 8\", \"356\u00f76=59, 2\"),\n    @(\"965\u00f79=107, 2\", \"981\u00f74=245, 1\"),\n    @(\"691\u00f74=172, 3\", \"515\u00f78=64, 3\"),\n    @(\"703\u00f75=140, 3\", \"159\u00f74=39, 3\"),\n    @(\"413\u00f76=68, 5\", \"420\u00f79=46, 6\"),\n    @(\"630\u00f76=105, 0\", \"795\u00f77=113, 4\"),\n    @(\"875\u00f75=175, 0\", \"640\u00f75=128, 0\"),\n    @(\"628\u00f75=125, 3\", \"367\u00f72=183, 1\"),\n    @(\"601\u00f74=150, 1\", \"561\u00f78=70, 1\"),\n    @(\"382\u00f78=47, 6\", \"554\u00f78=69, 2\"),\n    @(\"745\u00f78=93, 1\", \"602\u00f76=100, 2\"),\n    @(\"872\u00f74=218, 0\", \"544\u00f79=60, 4\"),\n    @(\"348\u00f77=49, 5\", \"250\u00f75=50, 0\"),\n    @(\"477\u00f78=59, 5\", \"546\u00f75=109, 1\"),\n    @(\"436\u00f73=145, 1\", \"729\u00f74=182, 1\"),\n    @(\"688\u00f72=344, 0\", \"744\u00f78=93, 0\"),\n    @(\"896\u00f74=224, 0\", \"841\u00f76=140, 1\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
